$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.450.34"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "2.426.75"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "512.43"
$ws.Range("E5").Value = "  -2.74%  "
$ws.Range("D6").Value = "129.39"
$ws.Range("E6").Value = "  -3.39%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "2.436.50"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").Value = "0.0951"
$ws.Range("E11").Value = "  -5.40%  "
$ws.Range("E12").Value = "  -3.66%  "
$ws.Range("D13").Value = "0.333"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").Value = "2.855.78"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "57.350.46"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "21.74"
$ws.Range("E16").Value = "  -3.10%  "
$ws.Range("E17").Value = "  -3.78%  "
$ws.Range("D18").Value = "2.432.01"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").Value = "315.07"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").Value = "63.61"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "0.405"
$ws.Range("E25").Value = "  -2.26%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("D29").Value = "168.73"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("E30").Value = "  -4.22%  "
$ws.Range("D31").Value = "6.23"
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("E37").Value = "  -5.16%  "
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("D39").Value = "36.09"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("D41").Value = "0.774"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").Value = "3.37"
$ws.Range("E42").Value = "  -4.65%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "267.64"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "4.89"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "119.75"
$ws.Range("E47").Value = "  -6.40%  "
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("D49").Value = "17.05"
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("D50").Value = "0.0210"
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("D51").Value = "16.47"
$ws.Range("E51").Value = "  -4.15%  "
